# Update the team-specific matrix (The Citadel_B) on Sheet1 with refreshed
# simulation results after adding more games / speeding up the simulate-game
# logic. Cell values below are the recomputed transition probabilities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2022471910112359
$ws.Cells.Item(2, 3).Value = 0.5168539325842697
$ws.Cells.Item(2, 10).Value = 0.03370786516853932
$ws.Cells.Item(2, 16).Value = 0.1685393258426966
$ws.Cells.Item(2, 19).Value = 0.07865168539325842
$ws.Cells.Item(3, 3).Value = 0.02173913043478261
$ws.Cells.Item(3, 10).Value = 0.0108695652173913
$ws.Cells.Item(3, 16).Value = 0.7826086956521739
$ws.Cells.Item(3, 19).Value = 0.1847826086956522
$ws.Cells.Item(4, 10).Value = 0.05714285714285714
$ws.Cells.Item(4, 16).Value = 0.6
$ws.Cells.Item(4, 19).Value = 0.3428571428571429
$ws.Cells.Item(6, 2).Value = 0.06310679611650485
$ws.Cells.Item(6, 4).Value = 0.01456310679611651
$ws.Cells.Item(6, 6).Value = 0.05825242718446602
$ws.Cells.Item(6, 10).Value = 0.2378640776699029
$ws.Cells.Item(6, 15).Value = 0.01941747572815534
$ws.Cells.Item(6, 17).Value = 0.1796116504854369
$ws.Cells.Item(6, 18).Value = 0.07281553398058252
$ws.Cells.Item(6, 19).Value = 0.354368932038835
$ws.Cells.Item(7, 2).Value = 0.0670391061452514
$ws.Cells.Item(7, 4).Value = 0.0111731843575419
$ws.Cells.Item(7, 6).Value = 0.0446927374301676
$ws.Cells.Item(7, 10).Value = 0.0893854748603352
$ws.Cells.Item(7, 15).Value = 0.03910614525139665
$ws.Cells.Item(7, 17).Value = 0.1899441340782123
$ws.Cells.Item(7, 18).Value = 0.07262569832402235
$ws.Cells.Item(7, 19).Value = 0.4860335195530726
$ws.Cells.Item(8, 2).Value = 0.06712962962962964
$ws.Cells.Item(8, 4).Value = 0.01388888888888889
$ws.Cells.Item(8, 5).Value = 0.002314814814814815
$ws.Cells.Item(8, 6).Value = 0.05092592592592592
$ws.Cells.Item(8, 10).Value = 0.1134259259259259
$ws.Cells.Item(8, 15).Value = 0.01388888888888889
$ws.Cells.Item(8, 17).Value = 0.2268518518518519
$ws.Cells.Item(8, 18).Value = 0.1203703703703704
$ws.Cells.Item(8, 19).Value = 0.3912037037037037
$ws.Cells.Item(9, 2).Value = 0.04487179487179487
$ws.Cells.Item(9, 6).Value = 0.108974358974359
$ws.Cells.Item(9, 10).Value = 0.1602564102564103
$ws.Cells.Item(9, 15).Value = 0.01282051282051282
$ws.Cells.Item(9, 17).Value = 0.1730769230769231
$ws.Cells.Item(9, 18).Value = 0.08974358974358974
$ws.Cells.Item(9, 19).Value = 0.4102564102564102
$ws.Cells.Item(10, 2).Value = 0.06747404844290658
$ws.Cells.Item(10, 4).Value = 0.02076124567474048
$ws.Cells.Item(10, 5).Value = 0.0008650519031141869
$ws.Cells.Item(10, 6).Value = 0.08044982698961937
$ws.Cells.Item(10, 10).Value = 0.1366782006920415
$ws.Cells.Item(10, 15).Value = 0.0259515570934256
$ws.Cells.Item(10, 17).Value = 0.2179930795847751
$ws.Cells.Item(10, 18).Value = 0.09688581314878893
$ws.Cells.Item(10, 19).Value = 0.3529411764705883
$ws.Cells.Item(11, 7).Value = 0.1691729323308271
$ws.Cells.Item(11, 10).Value = 0.08270676691729323
$ws.Cells.Item(11, 11).Value = 0.1917293233082707
$ws.Cells.Item(11, 12).Value = 0.5338345864661654
$ws.Cells.Item(11, 19).Value = 0.02255639097744361
$ws.Cells.Item(12, 7).Value = 0.7248322147651006
$ws.Cells.Item(12, 10).Value = 0.2080536912751678
$ws.Cells.Item(12, 11).Value = 0.01342281879194631
$ws.Cells.Item(12, 12).Value = 0.01342281879194631
$ws.Cells.Item(12, 19).Value = 0.04026845637583892
$ws.Cells.Item(13, 6).Value = 0.02380952380952381
$ws.Cells.Item(13, 7).Value = 0.7380952380952381
$ws.Cells.Item(13, 10).Value = 0.1666666666666667
$ws.Cells.Item(13, 19).Value = 0.07142857142857142
$ws.Cells.Item(15, 6).Value = 0.02487562189054726
$ws.Cells.Item(15, 8).Value = 0.1741293532338309
$ws.Cells.Item(15, 9).Value = 0.07960199004975124
$ws.Cells.Item(15, 10).Value = 0.3383084577114428
$ws.Cells.Item(15, 11).Value = 0.04975124378109453
$ws.Cells.Item(15, 15).Value = 0.03980099502487562
$ws.Cells.Item(15, 19).Value = 0.2935323383084577
$ws.Cells.Item(16, 6).Value = 0.008547008547008548
$ws.Cells.Item(16, 8).Value = 0.2051282051282051
$ws.Cells.Item(16, 9).Value = 0.08547008547008547
$ws.Cells.Item(16, 10).Value = 0.4188034188034188
$ws.Cells.Item(16, 11).Value = 0.07692307692307693
$ws.Cells.Item(16, 13).Value = 0.03418803418803419
$ws.Cells.Item(16, 15).Value = 0.05982905982905983
$ws.Cells.Item(16, 19).Value = 0.1111111111111111
$ws.Cells.Item(17, 6).Value = 0.006802721088435374
$ws.Cells.Item(17, 8).Value = 0.2018140589569161
$ws.Cells.Item(17, 9).Value = 0.06802721088435375
$ws.Cells.Item(17, 10).Value = 0.4421768707482993
$ws.Cells.Item(17, 11).Value = 0.07256235827664399
$ws.Cells.Item(17, 13).Value = 0.03174603174603174
$ws.Cells.Item(17, 15).Value = 0.07482993197278912
$ws.Cells.Item(17, 19).Value = 0.1020408163265306
$ws.Cells.Item(18, 6).Value = 0.01470588235294118
$ws.Cells.Item(18, 8).Value = 0.2009803921568628
$ws.Cells.Item(18, 9).Value = 0.07352941176470588
$ws.Cells.Item(18, 10).Value = 0.446078431372549
$ws.Cells.Item(18, 11).Value = 0.107843137254902
$ws.Cells.Item(18, 13).Value = 0.0196078431372549
$ws.Cells.Item(18, 15).Value = 0.04411764705882353
$ws.Cells.Item(18, 19).Value = 0.09313725490196079
$ws.Cells.Item(19, 6).Value = 0.01639344262295082
$ws.Cells.Item(19, 8).Value = 0.2194899817850638
$ws.Cells.Item(19, 9).Value = 0.07468123861566485
$ws.Cells.Item(19, 10).Value = 0.3743169398907104
$ws.Cells.Item(19, 11).Value = 0.1229508196721311
$ws.Cells.Item(19, 13).Value = 0.01730418943533698
$ws.Cells.Item(19, 15).Value = 0.06648451730418943
$ws.Cells.Item(19, 19).Value = 0.1083788706739526
